$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1858.5667
$ws.Range("I113").Value = 1794.5714
$ws.Range("J113").Value = 1878.0435
$ws.Range("K113").Value = 1794.5714
$ws.Range("L113").Value = 1878.0435
$ws.Range("M113").Value = 1459.4286
$ws.Range("N113").Value = -8386.0435

$ws.Range("H132").Value = 1921.4857
$ws.Range("I132").Value = 1919.1765
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5757.529500000001
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3227.529500000001
$ws.Range("N132").Value = -11060

$ws.Range("H135").Value = 1001.7778
$ws.Range("I135").Value = 568.5714
$ws.Range("K135").Value = 5117.1426
$ws.Range("M135").Value = -2582.1426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 286.3
$ws.Range("I5").Value = 151.57143
$ws.Range("K5").Value = 151.57143
$ws.Range("M5").Value = -39.57142999999999

$ws.Range("H61").Value = 2852.087
$ws.Range("I61").Value = 1409.4286
$ws.Range("K61").Value = 1409.4286
$ws.Range("M61").Value = -1197.4286

$ws.Range("H74").Value = 7221.1875
$ws.Range("I74").Value = 7221.1875
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 7221.1875
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -6347.1875

$ws.Range("H77").Value = 7221.1875
$ws.Range("I77").Value = 7221.1875
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 36105.9375
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -31737.9375

$ws.Range("H97").Value = 1201
$ws.Range("I97").Value = 987.2857
$ws.Range("J97").Value = 1699.6666
$ws.Range("K97").Value = 987.2857
$ws.Range("L97").Value = 1699.6666
$ws.Range("M97").Value = -491.2857
$ws.Range("N97").Value = -2691.6666

$ws.Range("H102").Value = 1568
$ws.Range("I102").Value = 1568
$ws.Range("K102").Value = 1568
$ws.Range("M102").Value = 54

$ws.Range("H122").Value = 1751.0834
$ws.Range("I122").Value = 1751.0834
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5253.2502
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -2803.2502

$ws.Range("H125").Value = 66379.164
$ws.Range("J125").Value = 66379.164
$ws.Range("L125").Value = 66379.164
$ws.Range("N125").Value = -76219.164

$ws.Range("H131").Value = 43698.75
$ws.Range("J131").Value = 43698.75
$ws.Range("L131").Value = 43698.75
$ws.Range("N131").Value = -53778.75

$ws.Range("H132").Value = 2485.4666
$ws.Range("I132").Value = 1652.2222
$ws.Range("J132").Value = 3735.3333
$ws.Range("K132").Value = 4956.6666
$ws.Range("L132").Value = 11205.9999
$ws.Range("M132").Value = -2426.6666
$ws.Range("N132").Value = -16265.9999

$ws.Range("H136").Value = 2852.087
$ws.Range("I136").Value = 1409.4286
$ws.Range("K136").Value = 4228.2858
$ws.Range("M136").Value = -1678.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 286.3
$ws.Range("I4").Value = 151.57143
$ws.Range("K4").Value = 151.57143
$ws.Range("M4").Value = -36.57142999999999

$ws.Range("H107").Value = 2518.7222
$ws.Range("I107").Value = 2540.1333
$ws.Range("J107").Value = 2411.6667
$ws.Range("K107").Value = 2540.1333
$ws.Range("L107").Value = 2411.6667
$ws.Range("M107").Value = -620.1333
$ws.Range("N107").Value = -6251.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 511.63416
$ws.Range("I107").Value = 433.35483
$ws.Range("J107").Value = 754.3
$ws.Range("K107").Value = 433.35483
$ws.Range("L107").Value = 754.3
$ws.Range("M107").Value = 1486.64517
$ws.Range("N107").Value = -4594.3

$ws.Range("H122").Value = 986.1177
$ws.Range("I122").Value = 985.25
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2955.75
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -505.75
$ws.Range("N122").Value = -7900

$ws.Range("H132").Value = 4533.8184
$ws.Range("I132").Value = 3486.5
$ws.Range("K132").Value = 10459.5
$ws.Range("M132").Value = -7929.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1833.3334
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 22500
$ws.Range("N132").Value = -27560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1915.5714
$ws.Range("I102").Value = 1915.5714
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1915.5714
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -293.5714

$ws.Range("H107").Value = 237.57895
$ws.Range("I107").Value = 233.94444
$ws.Range("J107").Value = 303
$ws.Range("K107").Value = 233.94444
$ws.Range("L107").Value = 303
$ws.Range("M107").Value = 1686.05556
$ws.Range("N107").Value = -4143

$ws.Range("H113").Value = 13268
$ws.Range("I113").Value = 933.3333
$ws.Range("J113").Value = 19435.334
$ws.Range("K113").Value = 933.3333
$ws.Range("L113").Value = 19435.334
$ws.Range("M113").Value = 1236.6667
$ws.Range("N113").Value = -23775.334

$ws.Range("H122").Value = 6251096
$ws.Range("I122").Value = 7143959.5
$ws.Range("J122").Value = 1050
$ws.Range("K122").Value = 21431878.5
$ws.Range("L122").Value = 3150
$ws.Range("M122").Value = -21429428.5
$ws.Range("N122").Value = -8050

$ws.Range("H132").Value = 7863.6875
$ws.Range("I132").Value = 8130.0713
$ws.Range("K132").Value = 24390.2139
$ws.Range("M132").Value = -21860.2139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 9195.111000000001
$ws.Range("I93").Value = 11322.214
$ws.Range("J93").Value = 1750.25
$ws.Range("K93").Value = 11322.214
$ws.Range("L93").Value = 1750.25
$ws.Range("M93").Value = -10074.214
$ws.Range("N93").Value = -4246.25

$ws.Range("H136").Value = 1854.6052
$ws.Range("I136").Value = 1423.3334
$ws.Range("K136").Value = 4270.0002
$ws.Range("M136").Value = -1720.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3480
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws.Range("H107").Value = 577
$ws.Range("I107").Value = 452.66666
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 1357.99998
$ws.Range("L107").Value = 2850
$ws.Range("M107").Value = 562.0000199999999
$ws.Range("N107").Value = -6690

$ws.Range("H122").Value = 2674.2778
$ws.Range("I122").Value = 2228.3572
$ws.Range("J122").Value = 4235
$ws.Range("K122").Value = 6685.071599999999
$ws.Range("L122").Value = 12705
$ws.Range("M122").Value = -4235.071599999999
$ws.Range("N122").Value = -17605

$ws.Range("H136").Value = 3826.6924
$ws.Range("I136").Value = 3865
$ws.Range("K136").Value = 11595
$ws.Range("M136").Value = -9045
